$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (Thurs, Nov 21): "Do Before Class" now links to the dask notebook
# instead of the old multi-link "Parallel ML with Dask / What is Dask? /
# OPTIONAL Full Dask Tutorial" block, and a new "In-Class Exercise" link
# is added pointing at the distributed-computing exercises page.
$ws.Range("C27").Formula = '''- `Distributed Computing with dask <distributed_computing.ipynb>`_'
$ws.Range("D27").Formula = '`Link <https://www.practicaldatascience.org/html/distributed_computing.html#Exercises>`_'

# Row 27 used to need a tall row to fit the long dask-tutorial text block;
# the new, shorter text only needs a couple of lines.
$ws.Rows.Item(27).RowHeight = 34

# The standalone "Regular Expressions" topic row at the bottom of the sheet
# is no longer needed (folded into row 26 above), so remove it entirely.
$ws.Rows.Item(31).Delete()

# Match the author's final selection/cursor position.
$ws.Range("D27").Select()
